# Apply the "456a3b4" data refresh to 南宁-漫展信息.xlsx:
#   - bump a couple of "want to go" counters
#   - insert a newly-scraped event ("原神x星铁x绝区零同人ONLY3.0") into the
#     "展览" (exhibitions) sheet and the "全部类型" (all types) sheet,
#     pushing the rows below it down by one.
#
# Helper: write a value into a cell while avoiding Excel's automatic
# "this looks like a date, let me reformat the cell" behaviour for strings
# like "2024-09-15". We flip the cell to Text, assign, then restore the
# cell to the default (unstyled) "Normal" style so the stored file matches
# a plain text cell exactly (no stray date format left behind).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Copy just the cell formatting (font/border/alignment) from $srcRange onto
# $dstRange without touching $dstRange's value, by round-tripping through
# the clipboard (Range.Style assignment resolves to the "Normal" named
# style here instead of an exact xf copy, which would drop the bold
# border style used for column A's row-number cells).
function Copy-CellFormat {
    param($srcRange, $dstRange)
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1): dimension A1:I4 -> A1:I5
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# "want to go" counters bumped
$ws1.Range("F2").Value = 1344
$ws1.Range("F3").Value = 2868

# The row that used to be row 4 (万圣漫控嘉年华10) moves down to row 5;
# its sequence number in column A increments from 3 to 4.
$ws1.Range("A5").Value = 4
Set-TextValue $ws1.Range("B5") "2024-11-02"
$ws1.Range("C5").Value = "南宁·万圣漫控嘉年华10"
$ws1.Range("D5").Value = "亭洪路45号 百益上河城"
$ws1.Range("E5").Value = "2024.11.02 11:00-11.03 22:00"
$ws1.Range("F5").Value = 258
$ws1.Range("G5").Value = 50
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
Set-TextValue $ws1.Range("I5") "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"
Copy-CellFormat $ws1.Range("A2") $ws1.Range("A5")

# New row 4: the newly scraped event.
$ws1.Range("A4").Value = 3
Set-TextValue $ws1.Range("B4") "2024-09-15"
$ws1.Range("C4").Value = "南宁·原神x星铁x绝区零同人ONLY3.0"
$ws1.Range("D4").Value = "亭洪路45号 百益上河城"
$ws1.Range("E4").Value = "2024.09.15 10:00-09.15 17:00"
$ws1.Range("F4").Value = 2
$ws1.Range("G4").Value = 60
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=90570"
Set-TextValue $ws1.Range("I4") "//i0.hdslb.com/bfs/openplatform/202408/sd7B5MV91723100089780.jpeg"
Copy-CellFormat $ws1.Range("A2") $ws1.Range("A4")

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4): dimension A1:I6 -> A1:I7
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# "want to go" counters bumped
$ws4.Range("F3").Value = 1344
$ws4.Range("F4").Value = 2868

# Old row 6 (万圣漫控嘉年华10) moves down to row 7, keeping its old
# sequence number (6) which now matches its new row position.
$ws4.Range("A7").Value = 6
Set-TextValue $ws4.Range("B7") "2024-11-02"
$ws4.Range("C7").Value = "南宁·万圣漫控嘉年华10"
$ws4.Range("D7").Value = "亭洪路45号 百益上河城"
$ws4.Range("E7").Value = "2024.11.02 11:00-11.03 22:00"
$ws4.Range("F7").Value = 258
$ws4.Range("G7").Value = 50
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
Set-TextValue $ws4.Range("I7") "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"
Copy-CellFormat $ws4.Range("A2") $ws4.Range("A7")

# Old row 5 (最后的莫西干人...) moves down to row 6, sequence number
# unchanged (5).
$ws4.Range("A6").Value = 5
Set-TextValue $ws4.Range("B6") "2024-10-04"
$ws4.Range("C6").Value = "南宁·《最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会》"
$ws4.Range("D6").Value = "福建园街道星光大道4号 南宁剧场"
$ws4.Range("E6").Value = "2024.10.04 20:00-10.04 21:30"
$ws4.Range("F6").Value = 6
$ws4.Range("G6").Value = 100
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=89039"
Set-TextValue $ws4.Range("I6") "//i0.hdslb.com/bfs/openplatform/202407/dudapgjU1720595605665.jpeg"
Copy-CellFormat $ws4.Range("A2") $ws4.Range("A6")

# Row 5 becomes the newly scraped event, sequence number unchanged (4).
$ws4.Range("A5").Value = 4
Set-TextValue $ws4.Range("B5") "2024-09-15"
$ws4.Range("C5").Value = "南宁·原神x星铁x绝区零同人ONLY3.0"
$ws4.Range("D5").Value = "亭洪路45号 百益上河城"
$ws4.Range("E5").Value = "2024.09.15 10:00-09.15 17:00"
$ws4.Range("F5").Value = 2
$ws4.Range("G5").Value = 60
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=90570"
Set-TextValue $ws4.Range("I5") "//i0.hdslb.com/bfs/openplatform/202408/sd7B5MV91723100089780.jpeg"
Copy-CellFormat $ws4.Range("A2") $ws4.Range("A5")

Write-Output "done"
